$d = $word.ActiveDocument

$d.Content.Find.Execute("85×73=", $true, $false, $false, $false, $false, $true, 1, $false, "63×57=", 2) | Out-Null
$d.Content.Find.Execute("44×25=", $true, $false, $false, $false, $false, $true, 1, $false, "64×81=", 2) | Out-Null
$d.Content.Find.Execute("98×71=", $true, $false, $false, $false, $false, $true, 1, $false, "21×57=", 2) | Out-Null
$d.Content.Find.Execute("16×91=", $true, $false, $false, $false, $false, $true, 1, $false, "60×59=", 2) | Out-Null
$d.Content.Find.Execute("69×65=", $true, $false, $false, $false, $false, $true, 1, $false, "86×96=", 2) | Out-Null
$d.Content.Find.Execute("82×68=", $true, $false, $false, $false, $false, $true, 1, $false, "53×37=", 2) | Out-Null
$d.Content.Find.Execute("30×71=", $true, $false, $false, $false, $false, $true, 1, $false, "97×67=", 2) | Out-Null
$d.Content.Find.Execute("47×49=", $true, $false, $false, $false, $false, $true, 1, $false, "67×67=", 2) | Out-Null
$d.Content.Find.Execute("87×14=", $true, $false, $false, $false, $false, $true, 1, $false, "71×60=", 2) | Out-Null
$d.Content.Find.Execute("67×53=", $true, $false, $false, $false, $false, $true, 1, $false, "11×28=", 2) | Out-Null
$d.Content.Find.Execute("31×43=", $true, $false, $false, $false, $false, $true, 1, $false, "96×18=", 2) | Out-Null
$d.Content.Find.Execute("90×61=", $true, $false, $false, $false, $false, $true, 1, $false, "48×34=", 2) | Out-Null
$d.Content.Find.Execute("24×13=", $true, $false, $false, $false, $false, $true, 1, $false, "91×25=", 2) | Out-Null
$d.Content.Find.Execute("19×37=", $true, $false, $false, $false, $false, $true, 1, $false, "47×68=", 2) | Out-Null
$d.Content.Find.Execute("62×74=", $true, $false, $false, $false, $false, $true, 1, $false, "45×48=", 2) | Out-Null
$d.Content.Find.Execute("58×19=", $true, $false, $false, $false, $false, $true, 1, $false, "90×92=", 2) | Out-Null
$d.Content.Find.Execute("58×67=", $true, $false, $false, $false, $false, $true, 1, $false, "52×25=", 2) | Out-Null
$d.Content.Find.Execute("12×59=", $true, $false, $false, $false, $false, $true, 1, $false, "31×99=", 2) | Out-Null
$d.Content.Find.Execute("58×27=", $true, $false, $false, $false, $false, $true, 1, $false, "68×84=", 2) | Out-Null
$d.Content.Find.Execute("16×37=", $true, $false, $false, $false, $false, $true, 1, $false, "21×58=", 2) | Out-Null
$d.Content.Find.Execute("90×30=", $true, $false, $false, $false, $false, $true, 1, $false, "37×47=", 2) | Out-Null
$d.Content.Find.Execute("71×70=", $true, $false, $false, $false, $false, $true, 1, $false, "66×11=", 2) | Out-Null
$d.Content.Find.Execute("81×79=", $true, $false, $false, $false, $false, $true, 1, $false, "72×42=", 2) | Out-Null
$d.Content.Find.Execute("23×92=", $true, $false, $false, $false, $false, $true, 1, $false, "21×20=", 2) | Out-Null
$d.Content.Find.Execute("55×60=", $true, $false, $false, $false, $false, $true, 1, $false, "37×90=", 2) | Out-Null
